$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1530.7391
$ws.Range("I33").Value = 1430.5883
$ws.Range("K33").Value = 1430.5883
$ws.Range("M33").Value = -1201.5883
$ws.Range("H107").Value = 389
$ws.Range("J107").Value = 550
$ws.Range("L107").Value = 550
$ws.Range("N107").Value = -4390
$ws.Range("H111").Value = 2507.5715
$ws.Range("I111").Value = 2004.6666
$ws.Range("J111").Value = 5525
$ws.Range("K111").Value = 6013.9998
$ws.Range("L111").Value = 16575
$ws.Range("M111").Value = -2946.9998
$ws.Range("N111").Value = -22709
$ws.Range("H112").Value = 1633.8889
$ws.Range("I112").Value = 533.3333
$ws.Range("K112").Value = 1599.9999
$ws.Range("M112").Value = -491.9999
$ws.Range("H113").Value = 45457216
$ws.Range("I113").Value = 14287694
$ws.Range("J113").Value = 100003870
$ws.Range("K113").Value = 14287694
$ws.Range("L113").Value = 100003870
$ws.Range("M113").Value = -14284440
$ws.Range("N113").Value = -100010378
$ws.Range("H137").Value = 6716
$ws.Range("I137").Value = 2834
$ws.Range("K137").Value = 8502
$ws.Range("M137").Value = -5952
$ws.Range("H138").Value = 3508.7354
$ws.Range("J138").Value = 3508.7354
$ws.Range("L138").Value = 10526.2062
$ws.Range("N138").Value = -20806.2062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2472.077
$ws.Range("I45").Value = 2070.7778
$ws.Range("J45").Value = 3375
$ws.Range("K45").Value = 2070.7778
$ws.Range("L45").Value = 3375
$ws.Range("M45").Value = -1693.7778
$ws.Range("N45").Value = -4129
$ws.Range("H61").Value = 12227928
$ws.Range("I61").Value = 14293239
$ws.Range("J61").Value = 180280
$ws.Range("K61").Value = 14293239
$ws.Range("L61").Value = 180280
$ws.Range("M61").Value = -14293027
$ws.Range("N61").Value = -180704
$ws.Range("H74").Value = 6585152.5
$ws.Range("I74").Value = 9260853
$ws.Range("J74").Value = 17523.727
$ws.Range("K74").Value = 9260853
$ws.Range("L74").Value = 17523.727
$ws.Range("M74").Value = -9259979
$ws.Range("N74").Value = -19271.727
$ws.Range("H77").Value = 6585152.5
$ws.Range("I77").Value = 9260853
$ws.Range("J77").Value = 17523.727
$ws.Range("K77").Value = 46304265
$ws.Range("L77").Value = 87618.63499999999
$ws.Range("M77").Value = -46299897
$ws.Range("N77").Value = -96354.63499999999
$ws.Range("H112").Value = 60290.75
$ws.Range("J112").Value = 60290.75
$ws.Range("L112").Value = 60290.75
$ws.Range("N112").Value = -63244.75
$ws.Range("H132").Value = 5940.9165
$ws.Range("I132").Value = 2325.125
$ws.Range("J132").Value = 13172.5
$ws.Range("K132").Value = 6975.375
$ws.Range("L132").Value = 39517.5
$ws.Range("M132").Value = -4445.375
$ws.Range("N132").Value = -44577.5
$ws.Range("H136").Value = 12227928
$ws.Range("I136").Value = 14293239
$ws.Range("J136").Value = 180280
$ws.Range("K136").Value = 42879717
$ws.Range("L136").Value = 540840
$ws.Range("M136").Value = -42877167
$ws.Range("N136").Value = -545940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 15000
$ws.Range("I5").Value = 15000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -14887
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 502956.38
$ws.Range("I31").Value = 5846.25
$ws.Range("K31").Value = 5846.25
$ws.Range("M31").Value = -5551.25
$ws.Range("H34").Value = 502956.38
$ws.Range("I34").Value = 5846.25
$ws.Range("K34").Value = 5846.25
$ws.Range("M34").Value = -5644.25
$ws.Range("H125").Value = 461388
$ws.Range("J125").Value = 461388
$ws.Range("L125").Value = 461388
$ws.Range("N125").Value = -466308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7120887.5
$ws.Range("J4").Value = 3236839
$ws.Range("L4").Value = 9710517
$ws.Range("N4").Value = -9710741
$ws.Range("H59").Value = 1709.091
$ws.Range("J59").Value = 2440
$ws.Range("L59").Value = 7320
$ws.Range("N59").Value = -8400
$ws.Range("H113").Value = 1557.0769
$ws.Range("J113").Value = 1787
$ws.Range("L113").Value = 5361
$ws.Range("N113").Value = -9701

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 101874
$ws.Range("J62").Value = 101874
$ws.Range("L62").Value = 101874
$ws.Range("N62").Value = -103246
$ws.Range("H65").Value = 101874
$ws.Range("J65").Value = 101874
$ws.Range("L65").Value = 305622
$ws.Range("N65").Value = -312486
$ws.Range("H122").Value = 1228.4286
$ws.Range("I122").Value = 981.55554
$ws.Range("J122").Value = 1672.8
$ws.Range("K122").Value = 2944.66662
$ws.Range("L122").Value = 5018.4
$ws.Range("M122").Value = -494.66662
$ws.Range("N122").Value = -9918.4
$ws.Range("H136").Value = 40081.5
$ws.Range("J136").Value = 40081.5
$ws.Range("L136").Value = 120244.5
$ws.Range("N136").Value = -125344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 147257.72
$ws.Range("I7").Value = 1594.6666
$ws.Range("J7").Value = 256505
$ws.Range("K7").Value = 1594.6666
$ws.Range("L7").Value = 256505
$ws.Range("M7").Value = -1482.6666
$ws.Range("N7").Value = -256729
$ws.Range("H16").Value = 1122.55
$ws.Range("I16").Value = 1115.25
$ws.Range("J16").Value = 1151.75
$ws.Range("K16").Value = 1115.25
$ws.Range("L16").Value = 1151.75
$ws.Range("M16").Value = -945.25
$ws.Range("N16").Value = -1491.75
$ws.Range("H63").Value = 57849.5
$ws.Range("J63").Value = 57849.5
$ws.Range("L63").Value = 57849.5
$ws.Range("N63").Value = -59347.5
$ws.Range("H66").Value = 57849.5
$ws.Range("J66").Value = 57849.5
$ws.Range("L66").Value = 173548.5
$ws.Range("N66").Value = -181036.5
$ws.Range("H68").Value = 1765
$ws.Range("I68").Value = 1765
$ws.Range("K68").Value = 1765
$ws.Range("M68").Value = -1016
$ws.Range("H71").Value = 1765
$ws.Range("I71").Value = 1765
$ws.Range("K71").Value = 8825
$ws.Range("M71").Value = -5081
$ws.Range("H126").Value = 147257.72
$ws.Range("I126").Value = 1594.6666
$ws.Range("J126").Value = 256505
$ws.Range("K126").Value = 4783.9998
$ws.Range("L126").Value = 769515
$ws.Range("M126").Value = -2313.9998
$ws.Range("N126").Value = -774455
$ws.Range("H132").Value = 45152.832
$ws.Range("I132").Value = 3882.6667
$ws.Range("K132").Value = 11648.0001
$ws.Range("M132").Value = -9118.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 113000
$ws.Range("J111").Value = 113000
$ws.Range("L111").Value = 113000
$ws.Range("N111").Value = -121180
$ws.Range("H113").Value = 1118.7059
$ws.Range("I113").Value = 1543.1
$ws.Range("K113").Value = 4629.299999999999
$ws.Range("M113").Value = -2459.299999999999
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H126").Value = 6850.8335
$ws.Range("I126").Value = 6401.3335
$ws.Range("K126").Value = 19204.0005
$ws.Range("M126").Value = -16734.0005
$ws.Range("H132").Value = 1609.7354
$ws.Range("I132").Value = 1457.44
$ws.Range("J132").Value = 2032.7778
$ws.Range("K132").Value = 4372.32
$ws.Range("L132").Value = 6098.3334
$ws.Range("M132").Value = -1842.32
$ws.Range("N132").Value = -11158.3334
$ws.Range("H136").Value = 648
$ws.Range("I136").Value = 648
$ws.Range("K136").Value = 1944
$ws.Range("M136").Value = 606